# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" for the
# 8b535931-459d-4105-b028-c620b0fdcf6c file (row 7) across all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-31 00:44:39"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-31 00:44:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-31 00:44:39"
